$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3113.2
$ws.Range("I76").Value = 3141.6667
$ws.Range("K76").Value = 3141.6667
$ws.Range("M76").Value = -2826.6667
$ws.Range("H79").Value = 3113.2
$ws.Range("I79").Value = 3141.6667
$ws.Range("K79").Value = 3141.6667
$ws.Range("M79").Value = -2049.6667
$ws.Range("H88").Value = 23186.586
$ws.Range("I88").Value = 48129.75
$ws.Range("J88").Value = 5579.647
$ws.Range("K88").Value = 48129.75
$ws.Range("L88").Value = 5579.647
$ws.Range("M88").Value = -47723.75
$ws.Range("N88").Value = -6391.647
$ws.Range("H91").Value = 23186.586
$ws.Range("I91").Value = 48129.75
$ws.Range("J91").Value = 5579.647
$ws.Range("K91").Value = 48129.75
$ws.Range("L91").Value = 5579.647
$ws.Range("M91").Value = -46725.75
$ws.Range("N91").Value = -8387.647000000001
$ws.Range("H116").Value = 3213.3333
$ws.Range("I116").Value = 3596
$ws.Range("J116").Value = 1300
$ws.Range("K116").Value = 3596
$ws.Range("L116").Value = 1300
$ws.Range("M116").Value = -154
$ws.Range("N116").Value = -8184
$ws.Range("H129").Value = 1107.7273
$ws.Range("I129").Value = 930
$ws.Range("J129").Value = 1125.5
$ws.Range("K129").Value = 2790
$ws.Range("L129").Value = 3376.5
$ws.Range("M129").Value = 2210
$ws.Range("N129").Value = -13376.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1735.8462
$ws.Range("I61").Value = 1426.6666
$ws.Range("K61").Value = 1426.6666
$ws.Range("M61").Value = -1214.6666
$ws.Range("H74").Value = 40692.96
$ws.Range("I74").Value = 84620.664
$ws.Range("J74").Value = 3040.6428
$ws.Range("K74").Value = 84620.664
$ws.Range("L74").Value = 3040.6428
$ws.Range("M74").Value = -83746.664
$ws.Range("N74").Value = -4788.6428
$ws.Range("H77").Value = 40692.96
$ws.Range("I77").Value = 84620.664
$ws.Range("J77").Value = 3040.6428
$ws.Range("K77").Value = 423103.32
$ws.Range("L77").Value = 15203.214
$ws.Range("M77").Value = -418735.32
$ws.Range("N77").Value = -23939.214
$ws.Range("H136").Value = 1735.8462
$ws.Range("I136").Value = 1426.6666
$ws.Range("K136").Value = 4279.9998
$ws.Range("M136").Value = -1729.9998
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2242.5
$ws.Range("I105").Value = 2156.6667
$ws.Range("K105").Value = 2156.6667
$ws.Range("M105").Value = -409.6667000000002
$ws.Range("H134").Value = 4137.355
$ws.Range("I134").Value = 3380.1924
$ws.Range("K134").Value = 10140.5772
$ws.Range("M134").Value = -7605.5772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2613.1694
$ws.Range("I31").Value = 1135.8966
$ws.Range("J31").Value = 4041.2
$ws.Range("K31").Value = 1135.8966
$ws.Range("L31").Value = 4041.2
$ws.Range("M31").Value = -840.8966
$ws.Range("N31").Value = -4631.2
$ws.Range("H34").Value = 2613.1694
$ws.Range("I34").Value = 1135.8966
$ws.Range("J34").Value = 4041.2
$ws.Range("K34").Value = 1135.8966
$ws.Range("L34").Value = 4041.2
$ws.Range("M34").Value = -933.8966
$ws.Range("N34").Value = -4445.2
$ws.Range("H141").Value = 49261.637
$ws.Range("J141").Value = 52187.8
$ws.Range("L141").Value = 52187.8
$ws.Range("N141").Value = -62547.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 966.55
$ws.Range("J131").Value = 992.70215
$ws.Range("L131").Value = 2978.10645
$ws.Range("N131").Value = -13058.10645
$ws.Range("H137").Value = 11449421
$ws.Range("I137").Value = 2888.5715
$ws.Range("J137").Value = 21465136
$ws.Range("K137").Value = 8665.7145
$ws.Range("L137").Value = 64395408
$ws.Range("M137").Value = -3565.7145
$ws.Range("N137").Value = -64405608
$ws.Range("H140").Value = 2551.9167
$ws.Range("I140").Value = 874.8182
$ws.Range("J140").Value = 21000
$ws.Range("K140").Value = 2624.4546
$ws.Range("L140").Value = 63000
$ws.Range("M140").Value = 2555.5454
$ws.Range("N140").Value = -73360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 15759.167
$ws.Range("J63").Value = 15759.167
$ws.Range("L63").Value = 15759.167
$ws.Range("N63").Value = -17131.167
$ws.Range("H66").Value = 15759.167
$ws.Range("J66").Value = 15759.167
$ws.Range("L66").Value = 47277.501
$ws.Range("N66").Value = -54141.501
$ws.Range("H80").Value = 2386.3076
$ws.Range("J80").Value = 2364.4
$ws.Range("L80").Value = 2364.4
$ws.Range("N80").Value = -4360.4
$ws.Range("H83").Value = 2386.3076
$ws.Range("J83").Value = 2364.4
$ws.Range("L83").Value = 11822
$ws.Range("N83").Value = -21806
$ws.Range("H126").Value = 16667720
$ws.Range("I126").Value = 19608872
$ws.Range("J126").Value = 1196.6666
$ws.Range("K126").Value = 58826616
$ws.Range("L126").Value = 3589.9998
$ws.Range("M126").Value = -58824146
$ws.Range("N126").Value = -8529.9998
$ws.Range("H132").Value = 4784.925
$ws.Range("I132").Value = 3083.2917
$ws.Range("J132").Value = 7337.375
$ws.Range("K132").Value = 9249.875100000001
$ws.Range("L132").Value = 22012.125
$ws.Range("M132").Value = -6719.875100000001
$ws.Range("N132").Value = -27072.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 676
$ws.Range("I22").Value = 448
$ws.Range("J22").Value = 736
$ws.Range("K22").Value = 448
$ws.Range("L22").Value = 736
$ws.Range("M22").Value = -153
$ws.Range("N22").Value = -1326
$ws.Range("H27").Value = 676
$ws.Range("I27").Value = 448
$ws.Range("J27").Value = 736
$ws.Range("K27").Value = 448
$ws.Range("L27").Value = 736
$ws.Range("M27").Value = -341
$ws.Range("N27").Value = -950
$ws.Range("H100").Value = 1890.0714
$ws.Range("I100").Value = 1572.2
$ws.Range("J100").Value = 2066.6667
$ws.Range("K100").Value = 1572.2
$ws.Range("L100").Value = 2066.6667
$ws.Range("M100").Value = -1031.2
$ws.Range("N100").Value = -3148.6667
$ws.Range("H123").Value = 20522.166
$ws.Range("J123").Value = 20522.166
$ws.Range("L123").Value = 20522.166
$ws.Range("N123").Value = -30322.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2528.7144
$ws.Range("I81").Value = 2466.8333
$ws.Range("J81").Value = 2900
$ws.Range("K81").Value = 4933.6666
$ws.Range("L81").Value = 5800
$ws.Range("M81").Value = -3872.6666
$ws.Range("N81").Value = -7922
$ws.Range("H84").Value = 2528.7144
$ws.Range("I84").Value = 2466.8333
$ws.Range("J84").Value = 2900
$ws.Range("K84").Value = 24668.333
$ws.Range("L84").Value = 29000
$ws.Range("M84").Value = -19364.333
$ws.Range("N84").Value = -39608
$ws.Range("H136").Value = 3749.5405
$ws.Range("I136").Value = 669.931
$ws.Range("J136").Value = 14913.125
$ws.Range("K136").Value = 2009.793
$ws.Range("L136").Value = 44739.375
$ws.Range("M136").Value = 540.2069999999999
$ws.Range("N136").Value = -49839.375

